# KIBON-2489 add property eigenleistungenGemeinde
#
# Inserts a new column into the "Data" sheet (between the "Delegationsmodell"
# block and the "Resultate" block, i.e. before the existing 4_Kantonsbeitrag
# column) carrying the new report field "4_eigenleistungen_Gemeinde" /
# "{eigenleistungenGemeinde}".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# The new column lands at sheet column 84 (CF); everything from the old
# column 84 onward (the "Resultate" group: 4_Kantonsbeitrag, ...) shifts one
# column to the right.
$newCol = 84

# Insert a new column before column 84; Excel shifts the existing "Resultate"
# columns (and their formatting) one column to the right automatically.
$ws.Columns.Item($newCol).Insert()

# The "Delegationsmodell" header was merged over CC6:CE6 (cols 81-83); extend
# it to also cover the freshly inserted column so it reads CC6:CF6.
$ws.Range($ws.Cells.Item(6, 81), $ws.Cells.Item(6, $newCol)).Merge()

# Match the new column's width to its neighbours in the Delegationsmodell
# block.
$ws.Columns.Item($newCol).ColumnWidth = $ws.Columns.Item($newCol - 1).ColumnWidth

# Populate the new header (row 7) / placeholder (row 8) cells.
$ws.Cells.Item(7, $newCol).Value = "4_eigenleistungen_Gemeinde"
$ws.Cells.Item(8, $newCol).Value = "{eigenleistungenGemeinde}"
